# Regenerate merged AHB files
# 1) Rename the "_old" / "_new" header-name suffixes to "_FV2210" / "_FV2304"
# 2) Turn the used range into a native Excel table ("Table1")
# 3) Freeze the header row (row 1)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headerRange = $ws.Range("A1:U1")
$stash = $ws.Range("A1000:U1000")

# --- 1. Rename header suffixes on the header row (row 1) ---
$headerRange.Replace("_old", "_FV2210")
$headerRange.Replace("_new", "_FV2304")

# --- Preserve the header row's existing look (bold/fill/border/wrap) ---
# around the table insertion: ListObjects.Add bakes whatever formatting is
# currently on the header into a dxf/table style, so stash it, reset to the
# workbook default, build the table, then restore the original formatting.
$headerRange.Copy()
$stash.PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false
$headerRange.Style = "Normal"

# --- 2. Convert the data range into a native Excel table ---
$dataRange = $ws.Range("A1:U63")
$lo = $ws.ListObjects.Add(1, $dataRange, $null, 1)
$lo.Name = "Table1"
$lo.TableStyle = ""

# --- restore the header formatting ---
$stash.Copy()
$headerRange.PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false
$stash.Clear()

# --- 3. Freeze panes above row 2 (i.e. freeze the header row) ---
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
